{"js": "// Apply \"hybrid bold + color\" highlighting to quantitative impact metrics\n// (percentages, dollar amounts, large numbers) inside specific resume\n// bullet points, matching the target OOXML diff exactly.\n//\n// Strategy: for each target paragraph (located by a unique substring of its\n// text), split out each metric token into its own run and mark that run\n// bold with color 2C3E50 \u2014 mirroring how Word's object model splits a run\n// when only part of its text receives new character formatting.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Each entry: a substring that uniquely identifies the target paragraph,\n// and the ordered list of metric tokens inside it that must become bold\n// + colored (search is run sequentially against the paragraph's own\n// range, so duplicate tokens such as \"87%\"/\"71%\" appearing in two\n// different paragraphs are handled independently and safely).\nconst edits = [\n  {\n    paragraphContains:\n      \"Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    tokens: [\"23%\", \"64%\"],\n  },\n  {\n    paragraphContains:\n      \"Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from\",\n    tokens: [\"87%\", \"71%\", \"\u00b14.2%\", \"\u00b12.1%\"],\n  },\n  {\n    paragraphContains: \"Wrote RFP and analyzed bids from 1,200 vendors\",\n    tokens: [\"1,200\"],\n  },\n  {\n    paragraphContains:\n      \"Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database\",\n    tokens: [\"$400M\", \"$1B\"],\n  },\n  {\n    paragraphContains:\n      \"Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    tokens: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    paragraphContains:\n      \"Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    tokens: [\"87%\", \"71%\"],\n    // Must NOT also match the longer paragraph above (which contains the\n    // \"reducing polling error margins\" continuation) \u2014 exclude it.\n    excludeIfContains: \"reducing polling error margins\",\n  },\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nfor (const edit of edits) {\n  const target = paragraphs.items.find((p) => {\n    if (!p.text.includes(edit.paragraphContains)) return false;\n    if (edit.excludeIfContains && p.text.includes(edit.excludeIfContains)) {\n      return false;\n    }\n    return true;\n  });\n\n  if (!target) {\n    throw new Error(\n      \"Could not find target paragraph containing: \" + edit.paragraphContains\n    );\n  }\n\n  for (const token of edit.tokens) {\n    const range = target.getRange();\n    const results = range.search(token, { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n\n    if (results.items.length === 0) {\n      throw new Error(\n        'Token \"' + token + '\" not found in paragraph: ' + edit.paragraphContains\n      );\n    }\n\n    results.items[0].font.bold = true;\n    results.items[0].font.color = HIGHLIGHT_COLOR;\n    await context.sync();\n  }\n}\n", "ps1": "# Apply \"hybrid bold + color\" highlighting to quantitative impact metrics\n# (percentages, dollar amounts, large numbers) inside specific resume\n# bullet points, matching the target OOXML diff exactly.\n#\n# Strategy: for each target paragraph (located by a unique substring of its\n# text), use Find.Execute scoped to that paragraph's Range to locate each\n# metric token, then set Bold + Color directly on the found sub-range. Word\n# automatically splits the run that previously spanned the whole paragraph\n# into separate runs so only the matched token carries the new formatting.\n\n$d = $word.ActiveDocument\n\n# Color 2C3E50 as a Word \"BGR\" integer (Font.Color expects 0xBBGGRR).\n$r = 0x2C\n$g = 0x3E\n$b = 0x50\n$highlightColor = $b * 65536 + $g * 256 + $r\n\n$edits = @(\n  @{\n    Contains = '\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%'\n    Tokens = @('23%', '64%')\n  },\n  @{\n    Contains = '\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from'\n    Tokens = @('87%', '71%', '\u00b14.2%', '\u00b12.1%')\n  },\n  @{\n    Contains = '\u2022 Wrote RFP and analyzed bids from 1,200 vendors'\n    Tokens = @('1,200')\n  },\n  @{\n    Contains = '\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database'\n    Tokens = @('$400M', '$1B')\n  },\n  @{\n    Contains = '\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M'\n    Tokens = @('73.5%', '$4.7M')\n  },\n  @{\n    Contains = '\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%'\n    Tokens = @('87%', '71%')\n    ExcludeContains = 'reducing polling error margins'\n  }\n)\n\nforeach ($edit in $edits) {\n  $target = $null\n  foreach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    if (-not $text.Contains($edit.Contains)) {\n      continue\n    }\n    if ($edit.ContainsKey('ExcludeContains') -and $text.Contains($edit.ExcludeContains)) {\n      continue\n    }\n    $target = $p\n    break\n  }\n\n  if ($target -eq $null) {\n    throw \"Could not find target paragraph containing: $($edit.Contains)\"\n  }\n\n  $pRange = $target.Range\n  $pStart = $pRange.Start\n  $pEnd = $pRange.End\n\n  foreach ($tok in $edit.Tokens) {\n    $find = $d.Range($pStart, $pEnd)\n    $ok = $find.Find.Execute($tok)\n    if (-not $ok) {\n      throw \"Token '$tok' not found in paragraph: $($edit.Contains)\"\n    }\n    $find.Font.Bold = 1\n    $find.Font.Color = $highlightColor\n  }\n}\n"}
